$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.441.95"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.918.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.90%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.009"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.55%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.008"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.53%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4824"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4071"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +2.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.013"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.33"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.918.43"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.066"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.248"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.62"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06875"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.86%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.010"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001040"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.82%  "
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("E20").Value = "  +0.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "29.428.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.650"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.74"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.180"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.155.02"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.64%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.688"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +10.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "155.82"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.111"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.69"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.016"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09607"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.649"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.554"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.374"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02282"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06100"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.58%  "
$ws.Range("E38").Value = "  +0.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.061"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.18%  "
$ws.Range("E40").Value = "  +6.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5973"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1846"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.285"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.402"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.27%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.46"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.13%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07601"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5596"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.952"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.68%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "118.15"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.429"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.18"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.01%  "
